$wb = $excel.ActiveWorkbook

# Sheet1: 原有设备一览
$ws1 = $wb.Worksheets.Item("原有设备一览")
$ws1.Range("B2").Value = "'1"
$ws1.Range("B3").Value = "'1"

# Sheet2: 原有设备能耗
$ws2 = $wb.Worksheets.Item("原有设备能耗")
$ws2.Range("B2").Value = "'1"
$ws2.Range("D2").Value = "'80644"
$ws2.Range("E2").Value = "'69589"
$ws2.Range("G2").Value = "(80644-69589)/80644*40%*132=7.238"
$ws2.Range("I2").Value = "'10.472"
$ws2.Range("J2").Value = "(132*1.15*0.8629+10.472)/(20.3*0.8629)=8.08"

$ws2.Range("B3").Value = "'1"
$ws2.Range("D3").Value = "'80644"
$ws2.Range("E3").Value = "'69589"
$ws2.Range("G3").Value = "(80644-69589)/80644*40%*132=7.238"
$ws2.Range("I3").Value = "'10.472"
$ws2.Range("J3").Value = "(132*1.15*0.8629+10.472)/(20.3*0.8629)=8.08"

# Sheet3: 能效对比
$ws3 = $wb.Worksheets.Item("能效对比")
$ws3.Range("B7").Value = "'5.9"
$ws3.Range("D7").Value = "'5.9"
$ws3.Range("B8").Value = "'0.0983"
$ws3.Range("D8").Value = "'0.0983"

$ws3.Range("B10").Value = 38.1868
$ws3.Range("C10").Value = 38.1868
$ws3.Range("D10").Value = 38.1868
$ws3.Range("E10").Value = 38.1868

$ws3.Range("B11").Value = 192461
$ws3.Range("C11").Value = 192461
$ws3.Range("D11").Value = 192461
$ws3.Range("E11").Value = 192461

$ws3.Range("B12").Value = 384922
$ws3.Range("C12").Value = 384922
$ws3.Range("D12").Value = 384922
$ws3.Range("E12").Value = 384922
